# "preparing simple data run and correcting plots"
# Update the "Coupling Parameters" sheet:
#  - max_permit_build_time (B11) becomes 4 (was 7)
#  - Look Ahead (B4) now derives from max_permit_build_time via formula (was a static 7)
#  - move the active selection to F9 (was H9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

$ws.Range("B11").Value = 4
$ws.Range("B4").Formula = "=B11"

$ws.Range("F9").Select()
